$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 271
$ws.Range("H2").Value = "living_rooms"
$ws.Range("I2").Value = "target"
$ws.Range("K2").Value = "j"
$ws.Range("L2").Value = "stimuli/img_vgh2g.png"
$ws.Range("M2").Value = 93.81395348837209
$ws.Range("N2").Value = 78.27906976744185
$ws.Range("O2").Value = 86.04651162790697
$ws.Range("P2").Value = 43
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = 10
$ws.Range("S2").Value = 10

# Row 3
$ws.Range("F3").Value = 272
$ws.Range("H3").Value = "living_rooms"
$ws.Range("I3").Value = "target"
$ws.Range("K3").Value = "j"
$ws.Range("L3").Value = "stimuli/img_j4ttn.png"
$ws.Range("M3").Value = 12.61904761904762
$ws.Range("N3").Value = 11.42857142857143
$ws.Range("O3").Value = 12.02380952380952
$ws.Range("P3").Value = 42
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 1

# Row 4
$ws.Range("F4").Value = 273
$ws.Range("H4").Value = "living_rooms"
$ws.Range("I4").Value = "target"
$ws.Range("K4").Value = "j"
$ws.Range("L4").Value = "stimuli/img_rg4in.png"
$ws.Range("M4").Value = 49.3695652173913
$ws.Range("N4").Value = 30.21739130434782
$ws.Range("O4").Value = 39.79347826086956
$ws.Range("P4").Value = 46
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = 3
$ws.Range("S4").Value = 3

# Row 5
$ws.Range("F5").Value = 274
$ws.Range("H5").Value = "living_rooms"
$ws.Range("I5").Value = "target"
$ws.Range("K5").Value = "j"
$ws.Range("L5").Value = "stimuli/img_9oofc.png"
$ws.Range("M5").Value = 82.47619047619048
$ws.Range("N5").Value = 65.5
$ws.Range("O5").Value = 73.98809523809524
$ws.Range("P5").Value = 42
$ws.Range("Q5").Value = 8
$ws.Range("R5").Value = 8
$ws.Range("S5").Value = 8

# Row 6
$ws.Range("F6").Value = 275
$ws.Range("H6").Value = "living_rooms"
$ws.Range("I6").Value = "target"
$ws.Range("K6").Value = "j"
$ws.Range("L6").Value = "stimuli/img_tn8ys.png"
$ws.Range("M6").Value = 86.70454545454545
$ws.Range("N6").Value = 72.4090909090909
$ws.Range("O6").Value = 79.55681818181819
$ws.Range("P6").Value = 44
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = 10
$ws.Range("S6").Value = 10

# Row 7
$ws.Range("F7").Value = 276
$ws.Range("H7").Value = "living_rooms"
$ws.Range("I7").Value = "target"
$ws.Range("K7").Value = "j"
$ws.Range("L7").Value = "stimuli/img_swq34.png"
$ws.Range("M7").Value = 64.11363636363636
$ws.Range("N7").Value = 43.04545454545455
$ws.Range("O7").Value = 53.57954545454545
$ws.Range("P7").Value = 44
$ws.Range("Q7").Value = 5
$ws.Range("R7").Value = 5
$ws.Range("S7").Value = 5

# Row 8
$ws.Range("F8").Value = 277
$ws.Range("H8").Value = "living_rooms"
$ws.Range("I8").Value = "target"
$ws.Range("K8").Value = "j"
$ws.Range("L8").Value = "stimuli/img_rych7.png"
$ws.Range("M8").Value = 30.4468085106383
$ws.Range("N8").Value = 23.4468085106383
$ws.Range("O8").Value = 26.9468085106383
$ws.Range("P8").Value = 47
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 2
$ws.Range("S8").Value = 2

# Row 9
$ws.Range("F9").Value = 278
$ws.Range("H9").Value = "bedrooms"
$ws.Range("I9").Value = "distractor"
$ws.Range("K9").Value = "f"
$ws.Range("L9").Value = "stimuli/img_3h4c9.png"
$ws.Range("M9").Value = 85.47619047619048
$ws.Range("N9").Value = 67.26190476190476
$ws.Range("O9").Value = 76.36904761904762
$ws.Range("P9").Value = 42
$ws.Range("Q9").Value = 9
$ws.Range("R9").Value = 9
$ws.Range("S9").Value = 9

# Row 10
$ws.Range("F10").Value = 279
$ws.Range("H10").Value = "living_rooms"
$ws.Range("I10").Value = "target"
$ws.Range("K10").Value = "j"
$ws.Range("L10").Value = "stimuli/img_2qhro.png"
$ws.Range("M10").Value = 81.73809523809524
$ws.Range("N10").Value = 62.73809523809524
$ws.Range("O10").Value = 72.23809523809524
$ws.Range("P10").Value = 42
$ws.Range("Q10").Value = 8
$ws.Range("R10").Value = 8
$ws.Range("S10").Value = 8

# Row 11
$ws.Range("F11").Value = 280
$ws.Range("H11").Value = "living_rooms"
$ws.Range("I11").Value = "target"
$ws.Range("K11").Value = "j"
$ws.Range("L11").Value = "stimuli/img_kq9s9.png"
$ws.Range("M11").Value = 62.30232558139535
$ws.Range("N11").Value = 39.97674418604651
$ws.Range("O11").Value = 51.13953488372093
$ws.Range("P11").Value = 43
$ws.Range("Q11").Value = 4
$ws.Range("R11").Value = 4
$ws.Range("S11").Value = 4

# Row 12
$ws.Range("F12").Value = 281
$ws.Range("H12").Value = "bedrooms"
$ws.Range("I12").Value = "distractor"
$ws.Range("K12").Value = "f"
$ws.Range("L12").Value = "stimuli/img_jp28n.png"
$ws.Range("M12").Value = 65.02564102564102
$ws.Range("N12").Value = 44.97435897435897
$ws.Range("O12").Value = 55
$ws.Range("P12").Value = 39
$ws.Range("Q12").Value = 4
$ws.Range("R12").Value = 4
$ws.Range("S12").Value = 4

# Row 13
$ws.Range("F13").Value = 282
$ws.Range("H13").Value = "bedrooms"
$ws.Range("I13").Value = "distractor"
$ws.Range("K13").Value = "f"
$ws.Range("L13").Value = "stimuli/img_d3t0o.png"
$ws.Range("M13").Value = 66.95121951219512
$ws.Range("N13").Value = 42.92682926829269
$ws.Range("O13").Value = 54.9390243902439
$ws.Range("P13").Value = 41
$ws.Range("Q13").Value = 4
$ws.Range("R13").Value = 4
$ws.Range("S13").Value = 4

# Row 14
$ws.Range("F14").Value = 283
$ws.Range("H14").Value = "bedrooms"
$ws.Range("I14").Value = "distractor"
$ws.Range("K14").Value = "f"
$ws.Range("L14").Value = "stimuli/img_twj5p.png"
$ws.Range("M14").Value = 67.71739130434783
$ws.Range("N14").Value = 42.08695652173913
$ws.Range("O14").Value = 54.90217391304348
$ws.Range("P14").Value = 46
$ws.Range("Q14").Value = 4
$ws.Range("R14").Value = 4
$ws.Range("S14").Value = 4

# Row 15
$ws.Range("F15").Value = 284
$ws.Range("H15").Value = "living_rooms"
$ws.Range("I15").Value = "target"
$ws.Range("K15").Value = "j"
$ws.Range("L15").Value = "stimuli/img_ra2nm.png"
$ws.Range("M15").Value = 70.75
$ws.Range("N15").Value = 50.375
$ws.Range("O15").Value = 60.5625
$ws.Range("P15").Value = 40
$ws.Range("Q15").Value = 6
$ws.Range("R15").Value = 6
$ws.Range("S15").Value = 6

# Row 16
$ws.Range("F16").Value = 285
$ws.Range("H16").Value = "living_rooms"
$ws.Range("I16").Value = "target"
$ws.Range("K16").Value = "j"
$ws.Range("L16").Value = "stimuli/img_ac0ey.png"
$ws.Range("M16").Value = 86.62222222222222
$ws.Range("N16").Value = 70.02222222222223
$ws.Range("O16").Value = 78.32222222222222
$ws.Range("P16").Value = 45
$ws.Range("Q16").Value = 9
$ws.Range("R16").Value = 9
$ws.Range("S16").Value = 9

# Row 17
$ws.Range("F17").Value = 286
$ws.Range("H17").Value = "living_rooms"
$ws.Range("I17").Value = "target"
$ws.Range("K17").Value = "j"
$ws.Range("L17").Value = "stimuli/img_zxvl3.png"
$ws.Range("M17").Value = 68.78260869565217
$ws.Range("N17").Value = 47.56521739130435
$ws.Range("O17").Value = 58.17391304347827
$ws.Range("P17").Value = 46
$ws.Range("Q17").Value = 5
$ws.Range("R17").Value = 5
$ws.Range("S17").Value = 5

# Row 18
$ws.Range("F18").Value = 287
$ws.Range("H18").Value = "living_rooms"
$ws.Range("I18").Value = "target"
$ws.Range("K18").Value = "j"
$ws.Range("L18").Value = "stimuli/img_24rt2.png"
$ws.Range("M18").Value = 55.26829268292683
$ws.Range("N18").Value = 34.19512195121951
$ws.Range("O18").Value = 44.73170731707317
$ws.Range("P18").Value = 41
$ws.Range("Q18").Value = 3
$ws.Range("R18").Value = 3
$ws.Range("S18").Value = 3

# Row 19
$ws.Range("F19").Value = 288
$ws.Range("H19").Value = "living_rooms"
$ws.Range("I19").Value = "target"
$ws.Range("K19").Value = "j"
$ws.Range("L19").Value = "stimuli/img_vh7v8.png"
$ws.Range("M19").Value = 78.70454545454545
$ws.Range("N19").Value = 59.63636363636363
$ws.Range("O19").Value = 69.17045454545455
$ws.Range("P19").Value = 44
$ws.Range("Q19").Value = 7
$ws.Range("R19").Value = 7
$ws.Range("S19").Value = 7

# Row 20
$ws.Range("F20").Value = 289
$ws.Range("H20").Value = "bedrooms"
$ws.Range("I20").Value = "distractor"
$ws.Range("K20").Value = "f"
$ws.Range("L20").Value = "stimuli/img_91csq.png"
$ws.Range("M20").Value = 50.44736842105263
$ws.Range("N20").Value = 28.34210526315789
$ws.Range("O20").Value = 39.39473684210526
$ws.Range("P20").Value = 38
$ws.Range("Q20").Value = 2
$ws.Range("R20").Value = 2
$ws.Range("S20").Value = 2

# Row 21
$ws.Range("F21").Value = 290
$ws.Range("H21").Value = "bedrooms"
$ws.Range("I21").Value = "distractor"
$ws.Range("K21").Value = "f"
$ws.Range("L21").Value = "stimuli/img_5il0t.png"
$ws.Range("M21").Value = 48.09523809523809
$ws.Range("N21").Value = 30.90476190476191
$ws.Range("O21").Value = 39.5
$ws.Range("P21").Value = 42
$ws.Range("Q21").Value = 2
$ws.Range("R21").Value = 2
$ws.Range("S21").Value = 2

# Row 22
$ws.Range("F22").Value = 291
$ws.Range("H22").Value = "living_rooms"
$ws.Range("I22").Value = "target"
$ws.Range("K22").Value = "j"
$ws.Range("L22").Value = "stimuli/img_syam3.png"
$ws.Range("M22").Value = 41.32432432432432
$ws.Range("N22").Value = 26.2972972972973
$ws.Range("O22").Value = 33.81081081081081
$ws.Range("P22").Value = 37
$ws.Range("Q22").Value = 2
$ws.Range("R22").Value = 2
$ws.Range("S22").Value = 2

# Row 23
$ws.Range("F23").Value = 292
$ws.Range("H23").Value = "living_rooms"
$ws.Range("I23").Value = "target"
$ws.Range("K23").Value = "j"
$ws.Range("L23").Value = "stimuli/img_iudc4.png"
$ws.Range("M23").Value = 73.625
$ws.Range("N23").Value = 52.275
$ws.Range("O23").Value = 62.95
$ws.Range("P23").Value = 40
$ws.Range("Q23").Value = 6
$ws.Range("R23").Value = 6
$ws.Range("S23").Value = 6

# Row 24
$ws.Range("F24").Value = 293
$ws.Range("H24").Value = "living_rooms"
$ws.Range("I24").Value = "target"
$ws.Range("K24").Value = "j"
$ws.Range("L24").Value = "stimuli/img_amsgw.png"
$ws.Range("M24").Value = 86.08510638297872
$ws.Range("N24").Value = 65.95744680851064
$ws.Range("O24").Value = 76.02127659574468
$ws.Range("P24").Value = 47
$ws.Range("Q24").Value = 9
$ws.Range("R24").Value = 9
$ws.Range("S24").Value = 9

# Row 25
$ws.Range("F25").Value = 294
$ws.Range("H25").Value = "living_rooms"
$ws.Range("I25").Value = "target"
$ws.Range("K25").Value = "j"
$ws.Range("L25").Value = "stimuli/img_xr3up.png"
$ws.Range("M25").Value = 76.24444444444444
$ws.Range("N25").Value = 55.88888888888889
$ws.Range("O25").Value = 66.06666666666666
$ws.Range("P25").Value = 45
$ws.Range("Q25").Value = 7
$ws.Range("R25").Value = 7
$ws.Range("S25").Value = 7

# Row 26
$ws.Range("F26").Value = 295
$ws.Range("H26").Value = "living_rooms"
$ws.Range("I26").Value = "target"
$ws.Range("K26").Value = "j"
$ws.Range("L26").Value = "stimuli/img_of8d6.png"
$ws.Range("M26").Value = 26.04878048780488
$ws.Range("N26").Value = 19.14634146341463
$ws.Range("O26").Value = 22.59756097560975
$ws.Range("P26").Value = 41
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = 1
$ws.Range("S26").Value = 1

# Row 27
$ws.Range("F27").Value = 296
$ws.Range("H27").Value = "living_rooms"
$ws.Range("I27").Value = "target"
$ws.Range("K27").Value = "j"
$ws.Range("L27").Value = "stimuli/img_rru0v.png"
$ws.Range("M27").Value = 56.45238095238095
$ws.Range("N27").Value = 39.42857142857143
$ws.Range("O27").Value = 47.94047619047619
$ws.Range("P27").Value = 42
$ws.Range("Q27").Value = 4
$ws.Range("R27").Value = 4
$ws.Range("S27").Value = 4
